$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Semester"
$ws.Range("B1").Value = "Year"
$ws.Range("C1").Value = "Class"
$ws.Range("D1").Value = "Department"
$ws.Range("E1").Value = "Size"
$ws.Range("F1").Value = "In_person"

$ws.Range("F1").Select()
